$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing rows 3..25 down to 4..26.
$ws.Rows.Item(3).Insert()

# Populate the freshly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44630
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104001
$ws.Range("J3").Value = "Granada"
$ws.Range("K3").Value = "Wonderfull"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 75
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 14 kilos empedrada"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1071
$ws.Range("T3").Value = 14
